# Applies the upload: resort Sheet1's results by dataset_name, and add a new
# "Sheet3" worksheet containing the Mean_imputed_* / train_data_control summary
# table (with the "accuracy" column left unformatted and "precision" shown
# with a Comma/thousands number style), matching the author's second upload.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- 1) Re-sort the existing Sheet1 results table by dataset_name (A) ascending ---
$sortRange = $ws1.Range("A1:F9")
$keyRange = $ws1.Range("A2:A9")
$sortRange.Sort($keyRange, 1)

# Tidy up the column widths (the data got wider after resaving in a newer Excel)
$ws1.Columns.Item(1).AutoFit()
$ws1.Columns.Item(2).AutoFit()

# Select the whole table, matching the saved selection state
$ws1.Range("A1:F9").Select()

# --- 2) Add the new worksheet ("Sheet3") right after Sheet1 ---
$ws3 = $wb.Worksheets.Add($null, $ws1)
$ws3.Name = "Sheet3"

# --- 3) Write the header + data rows ---
$headers = @("dataset_name","model","accuracy","precision","recall","f1")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws3.Cells.Item(1, $c + 1).Value2 = $headers[$c]
}

$rows = @(
    @("Mean_imputed_10","k-NearestNeighbour_classifier",0.76329999999999998,0.73,0.76,0.67),
    @("Mean_imputed_10","DecisionTree_classifier",0.76139999999999997,0.66,0.76,0.66),
    @("Mean_imputed_40","k-NearestNeighbour_classifier",0.76329999999999998,0.73,0.76,0.67),
    @("Mean_imputed_40","DecisionTree_classifier",0.74590000000000001,0.73,0.75,0.73),
    @("Mean_imputed_70","k-NearestNeighbour_classifier",0.7681,0.75,0.77,0.68),
    @("Mean_imputed_70","DecisionTree_classifier",0.74109999999999998,0.7,0.74,0.71),
    @("train_data_control","k-NearestNeighbour_classifier",0.76249999999999996,0.71,0.76,0.66),
    @("train_data_control","DecisionTree_classifier",0.76759999999999995,0.77,0.77,0.77)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $rows[$r]
    $excelRow = $r + 2
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws3.Cells.Item($excelRow, $c + 1).Value2 = $row[$c]
    }
}

# --- 4) Format column D (precision) with the built-in "Comma" style ---
$ws3.Range("D2:D9").Style = "Comma"

# --- 5) Column widths on the new sheet ---
$ws3.Columns.Item(1).ColumnWidth = 20.81640625
$ws3.Columns.Item(2).AutoFit()
$ws3.Columns.Item(3).AutoFit()
$ws3.Columns.Item(4).ColumnWidth = 9
$ws3.Columns.Item(5).AutoFit()

# --- 6) Make the new sheet the active tab / selection, matching the upload ---
$ws3.Range("C18").Select()
$ws3.Select()
